$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 18 (pushes old rows 18+ down by one; rows 22/23 become 23/24,
# and the corresponding merged cells shift automatically).
$ws.Rows.Item(18).Insert()

# The inserted row should look like the (old) row 17 both in style and content,
# so copy row 17 into the new row 18 first.
$ws.Range("B17:J17").Copy($ws.Range("B18:J18"))

# Row 17 should now take on row 16's look (same style as row 16), so copy
# row 16 into row 17.
$ws.Range("B16:J16").Copy($ws.Range("B17:J17"))

# Fix up the period labels: row16=2506 (was 2507), row17=2507 (now matches row16's
# old text automatically via the copy above), row18=2508 (was 2506, new period added).
$ws.Range("E16").Value = "2506"
$ws.Range("E18").Value = "2508"

# Update the "Valor Mora" total and "Cant. Periodos" count.
$ws.Range("E11").Value = 170820
$ws.Range("F13").Value = 3
